# "Tried to implement Penality Reward System (unfinished)"
#
# 1. Weekly Quantity sheet: remove the two weekly-order rows for the
#    PO weeks ending 2023-07-23 (45130.99999999999, qty 25) and
#    2023-07-30 (45137.99999999999, qty 15). Deleting the rows shifts
#    everything below up by two, so the sheet goes from 55 data+header
#    rows (A1:B55) down to 53 (A1:B53).
# 2. Monthly Trend sheet: the July 2023 month row (45138.99999999999)
#    requested-quantity total drops from 160 to 120 (reflecting the
#    removed weekly rows above).

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Rows("16:17").Delete()

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B6").Value = 120
